$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.012.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.47'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.25%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.522'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.22%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.93'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +11.55%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0614'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.29%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.869.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.640.67'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.45%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.67'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +26.39%  '

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.578'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.91'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.033.83'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.87'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.96'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.96%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.28%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.71'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.09%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.77'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.75'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.50%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.67'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.17%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0492'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.88%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.21%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.433.58'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.67'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.72%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.74%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.35%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.45%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.22'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.838'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.44%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.00'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.33%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '55.30'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.70%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.84%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.775.88'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.41'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.96%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.54%  '
